# SCO-T1-A07-A08.pptx : Aula T2-S01 e T2-S02, updates de codigo
#
# 1) Every content slide (2..22) has its main body "Rectangle" shape
#    (the one previously sitting at y=900000 or y=913284 EMU, x unchanged)
#    raised to y=841276 EMU (x stays the same).
# 2) On the slide containing the "SUBR1_R / FILL 4" assembly listing, the
#    three runs "SUBR1_R<TAB>" + "FILL 4" + "<TAB><TAB>" are merged back
#    into a single run "SUBR1_R<TAB>FILL 4<TAB><TAB>".

$p = $ppt.ActivePresentation

# Point value that, once the host engine rounds it to a single-precision
# float and converts back to EMU (Top/Left are EMU/12700 = points,
# stored as Single), reliably floors to exactly 841276 EMU.
$newTopPts = 66.242215

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        $topEmu = [math]::Round($sh.Top * 12700)
        if ($topEmu -eq 900000 -or $topEmu -eq 913284) {
            $sh.Top = $newTopPts
        }
    }
}

# Merge the split "SUBR1_R" / "FILL 4" / tabs runs back into one run.
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            $tr = $sh.TextFrame.TextRange
            $full = $tr.Text
            $idx = $full.IndexOf("SUBR1_R`tFILL 4`t`t")
            if ($idx -ge 0) {
                # Re-typing text on an auto-fit shape makes the host
                # re-layout/re-measure it, but the source text is unchanged
                # here (only run boundaries were merged, same characters) so
                # the rendered extent must come back out unchanged too
                # (original cy = 3968702 EMU).
                $sub = $tr.Characters($idx + 1, 16)
                $sub.Text = "SUBR1_R`tFILL 4`t`t"
                $sh.Height = 312.496258
            }
        }
    }
}
